$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2024-05-15 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-16 Thursday", 2)

# Update the division problems in the table, cell by cell, to avoid
# collisions between old/new values (e.g. 71/7= -> 80/5=, while
# 80/5= -> 92/2= is also a required change).
$tbl = $d.Tables.Item(1)

$values = @(
    @("18÷8=", "16÷9=", "10÷2=", "92÷2=", "98÷3="),
    @("30÷2=", "33÷2=", "39÷6=", "35÷9=", "36÷3="),
    @("51÷5=", "22÷4=", "91÷7=", "76÷3=", "77÷4="),
    @("77÷6=", "79÷2=", "69÷4=", "39÷3=", "72÷7="),
    @("81÷4=", "76÷3=", "37÷3=", "74÷8=", "80÷5=")
)

$rowIndexes = @(1, 5, 9, 13, 17)

for ($r = 0; $r -lt $rowIndexes.Length; $r++) {
    $tableRow = $rowIndexes[$r]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $tbl.Cell($tableRow, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $values[$r][$c - 1]
    }
}
